# Update Week 17 (row 3 = "R") target depth data for the Giants workbook.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 477   # Short Att
$wsOff.Range("C3").Value = 299   # Short Comp
$wsOff.Range("F3").Value = 12    # Short Int

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 488   # Short Att
$wsDef.Range("C3").Value = 363   # Short Comp
$wsDef.Range("D3").Value = 95    # Deep Att
$wsDef.Range("E3").Value = 40    # Deep Comp
$wsDef.Range("G3").Value = 3     # Deep Int
